$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: Table 2 (SND01), VOZ cell -> split "Masculina para el narrador y
# femenina para lo que dice el personaje" into 3 runs: "Masculina" + " 1" +
# " para el narrador y femenina para lo que dice el personaje"
# ---------------------------------------------------------------------------
$cell1 = $d.Tables.Item(2).Cell(2, 2)
$start1 = $cell1.Range.Start

$r1a = $d.Range($start1, $start1 + 9)
$r1a.InsertAfter(" 1")

# Force a run boundary between " 1" and the rest by toggling bold on/off
# (ends back at not-bold, matching the target formatting).
$r1b = $d.Range($start1 + 9, $start1 + 11)
$r1b.Font.Bold = 1
$r1b.Font.Bold = 0

# ---------------------------------------------------------------------------
# Hunk 2: Table 3 (SND02), VOZ cell -> "Masculina" becomes bold (paragraph
# mark + run) and gains a new bold+italic run " 2"
# ---------------------------------------------------------------------------
$cell2 = $d.Tables.Item(3).Cell(2, 2)
$cell2.Range.Paragraphs.Item(1).Range.Font.Bold = 1

$start2 = $d.Tables.Item(3).Cell(2, 2).Range.Start
$r2a = $d.Range($start2, $start2 + 9)
$r2a.InsertAfter(" 2")
$r2b = $d.Range($start2 + 9, $start2 + 11)
$r2b.Font.Bold = 0
$r2b.Font.Bold = 1

# ---------------------------------------------------------------------------
# Hunk 3: Table 4 (SND03), CÓDIGO DEL AUDIO cell -> wrap the whole paragraph
# content with a _GoBack bookmark (bookmarkStart before the first run,
# bookmarkEnd after the last run), preserving every existing run untouched.
# ---------------------------------------------------------------------------
$cell3 = $d.Tables.Item(4).Cell(1, 2)
$full3 = $cell3.Range
$r3 = $d.Range($full3.Start, $full3.End)
$xml3 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="60765BC0" w14:textId="3255C199" w:rsidR="002850E0" w:rsidRPr="002850E0" w:rsidRDefault="009C7336" w:rsidP="00EF5A78"><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r w:rsidRPr="002850E0"><w:rPr><w:i/></w:rPr><w:t>LE_</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>08</w:t></w:r><w:r w:rsidRPr="002850E0"><w:rPr><w:i/></w:rPr><w:t>_</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>05</w:t></w:r><w:r w:rsidRPr="002850E0"><w:rPr><w:i/></w:rPr><w:t>_CO_REC</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>10</w:t></w:r><w:r w:rsidRPr="002850E0"><w:rPr><w:i/></w:rPr><w:t>0</w:t></w:r><w:r w:rsidR="002850E0" w:rsidRPr="002850E0"><w:rPr><w:i/></w:rPr><w:t>_SND</w:t></w:r><w:r w:rsidR="002850E0"><w:rPr><w:i/></w:rPr><w:t>03</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r3.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Hunk 4: Table 4 (SND03), VOZ cell -> same treatment as hunk 2: "Masculina"
# becomes bold (paragraph mark + run) and gains a new bold+italic run " 2"
# ---------------------------------------------------------------------------
$cell4 = $d.Tables.Item(4).Cell(2, 2)
$cell4.Range.Paragraphs.Item(1).Range.Font.Bold = 1

$start4 = $d.Tables.Item(4).Cell(2, 2).Range.Start
$r4a = $d.Range($start4, $start4 + 9)
$r4a.InsertAfter(" 2")
$r4b = $d.Range($start4 + 9, $start4 + 11)
$r4b.Font.Bold = 0
$r4b.Font.Bold = 1

# ---------------------------------------------------------------------------
# Hunk 5: Table 6 (SND05), first paragraph of the dialogue cell -> drop the
# existing bookmarkStart/bookmarkEnd ("_GoBack") pair that trails the
# "Cova:] " run, leaving the run content itself untouched.
# ---------------------------------------------------------------------------
$cell5 = $d.Tables.Item(6).Cell(3, 1)
$full5 = $cell5.Range
$text5 = $full5.Text
$marker = "Cova:] "
$idx5 = $text5.IndexOf($marker)
$parEnd5 = $full5.Start + $idx5 + $marker.Length
$r5 = $d.Range($full5.Start, $parEnd5)
$xml5 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1206A39E" w14:textId="77777777" w:rsidR="002850E0" w:rsidRPr="00851DF5" w:rsidRDefault="002850E0" w:rsidP="002850E0"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:after="80" w:line="276" w:lineRule="auto"/><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r w:rsidRPr="00851DF5"><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">[Diálogo de Clemente Silva y Arturo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00851DF5"><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>Cova</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00851DF5"><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">:] </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r5.InsertXML($xml5)

"done"
